$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating "2022-Q3" (so it inherits
#    the exact same layout/header/styles), inserted right before it. This
#    yields the target sheet order: 总计, 2022-Q4, 2022-Q3, 2022-Q2.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("2022-Q3")
$ws3.Copy($ws3, $null)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Force the "numeric looking" text columns to be stored as text (matches the
# source data, which keeps these as strings) instead of being auto-converted
# to numbers. Column A (index) and H (rank) stay numeric; G18:G19 are the
# real numeric zeros further down so they are excluded here.
$newSheet.Range("B2:B19").NumberFormat = "@"
$newSheet.Range("D2:F19").NumberFormat = "@"
$newSheet.Range("G2:G17").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 2. Write the 2022-Q4 fund holdings data (rows 2-19).
# ---------------------------------------------------------------------------
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "501077"
$newSheet.Range("C2").Value = "富国创新企业灵活配置混合（LOF）A"
$newSheet.Range("D2").Value = "8.06"
$newSheet.Range("E2").Value = "88.91"
$newSheet.Range("F2").Value = "2.53"
$newSheet.Range("G2").Value = "0.2039"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "100029"
$newSheet.Range("C3").Value = "富国天成红利混合"
$newSheet.Range("D3").Value = "8.86"
$newSheet.Range("E3").Value = "73.21"
$newSheet.Range("F3").Value = "2.18"
$newSheet.Range("G3").Value = "0.1931"
$newSheet.Range("H3").Value = 10

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "673060"
$newSheet.Range("C4").Value = "西部利得景瑞灵活配置混合A"
$newSheet.Range("D4").Value = "3.09"
$newSheet.Range("E4").Value = "90.81"
$newSheet.Range("F4").Value = "4.74"
$newSheet.Range("G4").Value = "0.1465"
$newSheet.Range("H4").Value = 3

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "002376"
$newSheet.Range("C5").Value = "国寿安保核心产业灵活配置混合"
$newSheet.Range("D5").Value = "4.70"
$newSheet.Range("E5").Value = "88.55"
$newSheet.Range("F5").Value = "2.90"
$newSheet.Range("G5").Value = "0.1363"
$newSheet.Range("H5").Value = 9

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "014317"
$newSheet.Range("C6").Value = "广发价值领航一年持有混合A"
$newSheet.Range("D6").Value = "2.40"
$newSheet.Range("E6").Value = "93.40"
$newSheet.Range("F6").Value = "4.49"
$newSheet.Range("G6").Value = "0.1078"
$newSheet.Range("H6").Value = 9

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "009258"
$newSheet.Range("C7").Value = "西部利得景瑞灵活配置混合C"
$newSheet.Range("D7").Value = "1.73"
$newSheet.Range("E7").Value = "90.81"
$newSheet.Range("F7").Value = "4.74"
$newSheet.Range("G7").Value = "0.0820"
$newSheet.Range("H7").Value = 3

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "016588"
$newSheet.Range("C8").Value = "富国融甄混合A"
$newSheet.Range("D8").Value = "4.57"
$newSheet.Range("E8").Value = "29.63"
$newSheet.Range("F8").Value = "1.49"
$newSheet.Range("G8").Value = "0.0681"
$newSheet.Range("H8").Value = 9

$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "014318"
$newSheet.Range("C9").Value = "广发价值领航一年持有混合C"
$newSheet.Range("D9").Value = "0.65"
$newSheet.Range("E9").Value = "93.40"
$newSheet.Range("F9").Value = "4.49"
$newSheet.Range("G9").Value = "0.0292"
$newSheet.Range("H9").Value = 9

$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "016589"
$newSheet.Range("C10").Value = "富国融甄混合C"
$newSheet.Range("D10").Value = "1.96"
$newSheet.Range("E10").Value = "29.63"
$newSheet.Range("F10").Value = "1.49"
$newSheet.Range("G10").Value = "0.0292"
$newSheet.Range("H10").Value = 9

$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "004258"
$newSheet.Range("C11").Value = "国寿安保稳嘉混合A"
$newSheet.Range("D11").Value = "2.16"
$newSheet.Range("E11").Value = "23.32"
$newSheet.Range("F11").Value = "1.33"
$newSheet.Range("G11").Value = "0.0287"
$newSheet.Range("H11").Value = 3

$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "004301"
$newSheet.Range("C12").Value = "国寿安保稳信混合A"
$newSheet.Range("D12").Value = "1.48"
$newSheet.Range("E12").Value = "22.04"
$newSheet.Range("F12").Value = "1.73"
$newSheet.Range("G12").Value = "0.0256"
$newSheet.Range("H12").Value = 2

$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "005175"
$newSheet.Range("C13").Value = "国寿安保消费新蓝海灵活配置混合"
$newSheet.Range("D13").Value = "0.70"
$newSheet.Range("E13").Value = "89.58"
$newSheet.Range("F13").Value = "3.16"
$newSheet.Range("G13").Value = "0.0221"
$newSheet.Range("H13").Value = 6

$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "015694"
$newSheet.Range("C14").Value = "瑞达策略优选混合A"
$newSheet.Range("D14").Value = "0.08"
$newSheet.Range("E14").Value = "76.83"
$newSheet.Range("F14").Value = "2.53"
$newSheet.Range("G14").Value = "0.0020"
$newSheet.Range("H14").Value = 8

$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "015849"
$newSheet.Range("C15").Value = "富国创新企业灵活配置混合（LOF）C"
$newSheet.Range("D15").Value = "0.01"
$newSheet.Range("E15").Value = "88.91"
$newSheet.Range("F15").Value = "2.53"
$newSheet.Range("G15").Value = "0.0003"
$newSheet.Range("H15").Value = 10

$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "004302"
$newSheet.Range("C16").Value = "国寿安保稳信混合C"
$newSheet.Range("D16").Value = "0.01"
$newSheet.Range("E16").Value = "22.04"
$newSheet.Range("F16").Value = "1.73"
$newSheet.Range("G16").Value = "0.0002"
$newSheet.Range("H16").Value = 2

$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = "004259"
$newSheet.Range("C17").Value = "国寿安保稳嘉混合C"
$newSheet.Range("D17").Value = "0.01"
$newSheet.Range("E17").Value = "23.32"
$newSheet.Range("F17").Value = "1.33"
$newSheet.Range("G17").Value = "0.0001"
$newSheet.Range("H17").Value = 3

$newSheet.Range("A18").Value = 16
$newSheet.Range("B18").Value = "015406"
$newSheet.Range("C18").Value = "国寿安保稳信混合E"
$newSheet.Range("D18").Value = "0.00"
$newSheet.Range("E18").Value = "22.04"
$newSheet.Range("F18").Value = "1.73"
$newSheet.Range("G18").Value = 0
$newSheet.Range("H18").Value = 2

$newSheet.Range("A19").Value = 17
$newSheet.Range("B19").Value = "015695"
$newSheet.Range("C19").Value = "瑞达策略优选混合C"
$newSheet.Range("D19").Value = "0.00"
$newSheet.Range("E19").Value = "76.83"
$newSheet.Range("F19").Value = "2.53"
$newSheet.Range("G19").Value = 0
$newSheet.Range("H19").Value = 8

# ---------------------------------------------------------------------------
# 3. Rows 2-5 of column A already carry the bold/centred/bordered "index"
#    style inherited from the "2022-Q3" copy. Rows 6-19 are brand new, so
#    apply the same look (bold, centered, top-aligned, thin border) to them.
# ---------------------------------------------------------------------------
$rngA = $newSheet.Range("A6:A19")
$rngA.Font.Bold = $true
$rngA.HorizontalAlignment = -4108
$rngA.VerticalAlignment = -4160
$rngA.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 4. Update the "总计" (summary) sheet: a new first row for 2022-Q4, and the
#    previous two rows shift down by one.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 18
$wsTotal.Range("D2").Value = 1.08

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 4
$wsTotal.Range("D3").Value = 0.32

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q2"
$wsTotal.Range("C4").Value = 12
$wsTotal.Range("D4").Value = 0.64

$wsTotal.Range("A4").Font.Bold = $true
$wsTotal.Range("A4").HorizontalAlignment = -4108
$wsTotal.Range("A4").VerticalAlignment = -4160
$wsTotal.Range("A4").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 5. Restore the originally-active tab: "2022-Q2" (now the last sheet).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
